{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n// \"\u00a9 2020 . Contact: ...\" footer block (and the blank paragraph that\n// precedes it) that used to follow the \"Requisitos\" section of the page.\n//\n// Anchor on the last real line of the \"Requisitos\" section\n// (\"LOB1004: C\u00e1lculo II (Requisito fraco)\") and walk forward through the\n// three paragraphs that must go away:\n//   1. an empty paragraph\n//   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\n// A further empty paragraph (and the page-break paragraph after it) stay\n// untouched.\n\nconst anchorResults = context.document.body.search(\n  \"LOB1004: C\u00e1lculo II (Requisito fraco)\",\n  { matchCase: true }\n);\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length > 0) {\n  const anchorParagraph = anchorResults.items[0].paragraphs.getFirst();\n\n  const blankParagraph = anchorParagraph.getNext();\n  const jupiterParagraph = blankParagraph.getNext();\n  const copyrightParagraph = jupiterParagraph.getNext();\n\n  jupiterParagraph.load(\"text\");\n  copyrightParagraph.load(\"text\");\n  await context.sync();\n\n  // Sanity-check before deleting, so we never nuke the wrong paragraphs.\n  if (\n    jupiterParagraph.text === \"Ver no Jupiter Salvar em pdf Salvar em docx\" &&\n    copyrightParagraph.text.indexOf(\"Powered by Jekyll and Github pages\") !== -1\n  ) {\n    copyrightParagraph.delete();\n    jupiterParagraph.delete();\n    blankParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" /\n# \"\u00a9 2020 . Contact: ...\" footer block (and the blank paragraph that\n# precedes it) that used to follow the \"Requisitos\" section of the page.\n#\n# Anchor on the last real line of the \"Requisitos\" section\n# (\"LOB1004: C\u00e1lculo II (Requisito fraco)\") and walk forward through the\n# three paragraphs that must go away:\n#   1. an empty paragraph\n#   2. \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3. \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n# A further empty paragraph (and the page-break paragraph after it) stay\n# untouched.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"LOB1004: C\u00e1lculo II (Requisito fraco)\")\n\nif ($found) {\n    $anchorParagraph = $searchRange.Paragraphs(1)\n\n    $blankParagraph = $anchorParagraph.Next()\n    $jupiterParagraph = $blankParagraph.Next()\n    $copyrightParagraph = $jupiterParagraph.Next()\n\n    $jupiterText = $jupiterParagraph.Range.Text.TrimEnd([char]13, [char]7)\n    $copyrightText = $copyrightParagraph.Range.Text.TrimEnd([char]13, [char]7)\n\n    # Sanity-check before deleting, so we never nuke the wrong paragraphs.\n    if ($jupiterText -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\" -and\n        $copyrightText -like \"*Powered by Jekyll and Github pages*\") {\n        $copyrightParagraph.Range.Delete()\n        $jupiterParagraph.Range.Delete()\n        $blankParagraph.Range.Delete()\n    }\n}\n"}
